$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $target) {
        $cell.Value = $replacement
    }
}
